$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old dbExcel/Neo4jData column),
# shifting the old B -> C and old C -> D.
$ws.Range("B1").EntireColumn.Insert()

# --- Header row ---
$ws.Range("B1").Value = "StatQuery"

# --- Data row ---
$statQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Staffordshire Bull Terrier']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

$ws.Range("B2").Value = $statQuery

# Match the wrap-text style used by the existing query cell (A2).
$ws.Range("B2").WrapText = $true

# --- Column widths ---
# Columns A, C and D already carry their original widths forward
# automatically (the insert shifts them, widths untouched). Only the
# brand-new column B needs an explicit width, matching column A's
# (75.81640625). The host's ColumnWidth setter quantizes to whole
# pixels before it re-derives the stored "characters" width, so we
# dial in the input that lands closest to the true target after that
# round-trip.
$ws.Range("B1").EntireColumn.ColumnWidth = 74.98307291666667

# The sheet view no longer scrolls to B2 / keeps B2 selected.
$ws.Range("B2").Select()
